$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "item_code"
$ws.Range("B1").Value = "item_name"

$ws.Columns.Item(1).ColumnWidth = 11.83

$ws.Range("B2").Select() | Out-Null
